$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "3f84269c-689e-4f90-92f4-0a0cb2e2db55.md"
$wsOverview.Range("B2").Value = "e2e\3f84269c-689e-4f90-92f4-0a0cb2e2db55.md"
$wsOverview.Range("G2").Value = "2016-10-19 11:25:56"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("A2").Value = "3f84269c-689e-4f90-92f4-0a0cb2e2db55.md"
$wsZhCn.Range("G2").Value = "3f84269c-689e-4f90-92f4-0a0cb2e2db55.3fca14f3ad45cc1f67541aa7b604ac828631b9cd.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-10-19 11:25:45"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("A2").Value = "3f84269c-689e-4f90-92f4-0a0cb2e2db55.md"
$wsDeDe.Range("G2").Value = "3f84269c-689e-4f90-92f4-0a0cb2e2db55.3fca14f3ad45cc1f67541aa7b604ac828631b9cd.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-10-19 11:25:56"
